$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width changes ---
# Column E (5): 14 -> 36
$ws.Columns.Item(5).ColumnWidth = 35.1667
# Column T (20): 32 -> 30
$ws.Columns.Item(20).ColumnWidth = 29.1667

# --- Row 2 data changes ---

# A2: TST -> AA1 (plain text)
$ws.Range("A2").Value = "AA1"

# B2: 08/04/25 -> 08/07/25 (text that looks like a date; force text so it is not
# auto-converted into a real Excel date serial number)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "08/07/25"
$ws.Range("B2").Style = "Normal"

# C2: testraj -> 2025-08-14 (text that looks like a date; force text)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2025-08-14"
$ws.Range("C2").Style = "Normal"

# D2: CAPAIR -> NORWRI
$ws.Range("D2").Value = "NORWRI"

# E2: Captive Aire -> Norman S. Wright Mech. Equip. LLC.
$ws.Range("E2").Value = "Norman S. Wright Mech. Equip. LLC."

# F2 unchanged ("I")

# G2: 2457243 -> 127310 (text that looks numeric; force text)
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "127310"
$ws.Range("G2").Style = "Normal"

# H2: 08/04/25 -> 08/07/25 (text date; force text)
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "08/07/25"
$ws.Range("H2").Style = "Normal"

# I2: 48529.29 -> 74143.65 (text that looks numeric; force text)
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "74143.65"
$ws.Range("I2").Style = "Normal"

# J2: 4416.85 -> 6663.65 (text that looks numeric; force text)
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "6663.65"
$ws.Range("J2").Style = "Normal"

# K2: 1021.25 -> removed entirely (cell no longer present)
$ws.Range("K2").ClearContents()

# L2: 43091.19 -> 67480.00 (text that looks numeric; force text, keep trailing zeros)
$ws.Range("L2").NumberFormat = "@"
$ws.Range("L2").Value = "67480.00"
$ws.Range("L2").Style = "Normal"

# M2: 1412 -> 1410 (real number)
$ws.Range("M2").Value = 1410

# P2: 5030 -> 5040 (text that looks numeric; force text)
$ws.Range("P2").NumberFormat = "@"
$ws.Range("P2").Value = "5040"
$ws.Range("P2").Style = "Normal"

# Q2: 320 -> 330 (real number)
$ws.Range("Q2").Value = 330

# R2: M -> E
$ws.Range("R2").Value = "E"

# T2: captive aire_1754943103394.pdf -> 127310 nsw_1755208409430.pdf
$ws.Range("T2").Value = "127310 nsw_1755208409430.pdf"
